# Jun's file updates for all IO data and others
#
# The "ISIC 20T21" column on the EVCRSbRIC sheet is split into two
# separate columns: "ISIC 20" and "ISIC 21". This is done by inserting a
# new column immediately before the existing "ISIC 20T21" column (which
# shifts it one column to the right), then labeling the two resulting
# header cells and giving the newly inserted data cell the same default
# value (0) as its neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVCRSbRIC")

# Column K currently holds "ISIC 20T21". Insert a new blank column there,
# pushing the old column (and everything after it) one place to the right.
$ws.Columns("K").Insert()

# Label the new column and relabel the (now shifted) old column.
$ws.Range("K1").Value = "ISIC 20"
$ws.Range("L1").Value = "ISIC 21"

# The data row mirrors every other ISIC column with a default of 0.
$ws.Range("K2").Value = 0
